# The source report has several duplicate-product rows (rows that share the
# same item-code / description in column C, from a stock-movement report
# that got appended to out of order). The fix re-aligns the per-row
# "Closing/Rate/Qty/Value" figures (columns B, E, F, G) so each duplicate
# row carries the figures that actually belong to it: the values found in
# each block of duplicate rows are cyclically rotated by one position
# (each row takes the B/E/F/G figures that were sitting one row above it,
# wrapping from the last row of the block back to the first).
#
# Columns A (serial no), C (item name), D (rate) are untouched - they are
# correct/identical across the duplicate rows already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a contiguous block of duplicate-item rows (in sheet
# order) whose B/E/F/G values need to be rotated by one position.
$blocks = @(
    @(136,137),
    @(163,164),
    @(246,247),
    @(277,278),
    @(292,293),
    @(294,295,296),
    @(299,300),
    @(311,312),
    @(356,357),
    @(467,468),
    @(472,473),
    @(479,480),
    @(485,486),
    @(732,733)
)

$cols = @("B","E","F","G")

foreach ($block in $blocks) {
    $n = $block.Count

    # Snapshot the current B/E/F/G values for every row in this block
    # before writing anything (so later writes don't clobber values we
    # still need to read for later rows).
    $oldVals = @{}
    for ($i = 0; $i -lt $n; $i++) {
        $r = $block[$i]
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
        }
        $oldVals[$i] = $rowVals
    }

    # Rotate: row i receives the old values of row (i-1), wrapping around
    # so row 0 receives the old values of the last row in the block.
    for ($i = 0; $i -lt $n; $i++) {
        $r = $block[$i]
        $srcIdx = ($i - 1 + $n) % $n
        $src = $oldVals[$srcIdx]
        foreach ($c in $cols) {
            $ws.Cells.Item($r, $c).Value = $src[$c]
        }
    }
}
